# Tester Bulk Upload Excel format - "add upload option" edit
#
# Summary of the change being applied:
#  - Remove the "Username*" column (old column N) entirely - accounts are now
#    created automatically, so the bulk-upload sheet no longer collects a
#    username from the uploader.
#  - A handful of header labels toggle whether they are "required" (marked
#    with a trailing "*" and rendered in red) or "optional" (no "*",
#    rendered in black).
#  - The rich-text "Prefered Contact Method* (Phone or Email)" header is
#    simplified to a single plain-text run "Prefered Contact Method (Phone
#    or Email)".
#  - "Professional Registration No\n(if available)" gains a "*".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Useful colour constants (OLE BGR-packed decimal, matches Font.Color).
$RED   = 1974729   # RGB(C9,21,1E) - required field colour
$BLACK = 0          # RGB(00,00,00) - optional field colour
$xlCenter = -4108

# ------------------------------------------------------------------
# 1. Drop the "Username*" column (old column N). Everything to its
#    right (Password.. Facility In charge Email) shifts one column left.
# ------------------------------------------------------------------
$ws.Columns("N").Delete()

# ------------------------------------------------------------------
# 2. Update header text - toggle "*" markers / simplify rich text.
#    (Columns below use the POST-delete lettering, i.e. what used to be
#    column O is now column N, etc.)
# ------------------------------------------------------------------
$ws.Range("A1").Value2 = "Professional Registration No*`n(if available)"
$ws.Range("D1").Value2 = "Last Name (Surname)"
$ws.Range("G1").Value2 = "Type of HIV Test Modality/Point"
$ws.Range("J1").Value2 = "Prefered Contact Method (Phone or Email)"
$ws.Range("K1").Value2 = "Current Job Title"
$ws.Range("M1").Value2 = "Time Worked As Tester"
$ws.Range("N1").Value2 = "Password"
$ws.Range("O1").Value2 = "Testing Site In charge Name"
$ws.Range("P1").Value2 = "Testing Site In charge Phone"
$ws.Range("R1").Value2 = "Facility In charge Name"
$ws.Range("S1").Value2 = "Facility In charge Phone"

# ------------------------------------------------------------------
# 3. Re-colour headers: required (has "*") => red, optional => black.
#    Most cells only need a font-colour change (wrap/centre already
#    correct), which is applied directly.
# ------------------------------------------------------------------

# Now-required field (gained a "*") -> red
$ws.Range("A1").Font.Color = $RED

# Now-optional fields (lost their "*") -> black
foreach ($addr in @("C1","D1","G1","J1","K1","M1","N1","O1","P1","R1","S1")) {
    $ws.Range($addr).Font.Color = $BLACK
}

# D1 ("Last Name (Surname)") previously rendered left-aligned / not
# wrapped; bring it in line with the other black headers (centered,
# wrapped) to match its new column neighbours.
$ws.Range("D1").HorizontalAlignment = $xlCenter
$ws.Range("D1").WrapText = $true

# H1 ("Phone*") stays red/required but switches off word-wrap. The
# WrapText setter in this runtime only registers a change when going
# false -> true, so flip it on then build a dedicated no-wrap style and
# apply it, then restore the centering that a raw style application
# doesn't carry across.
$noWrapStyleName = "PhoneHeaderNoWrap"
$phoneStyle = $wb.Styles.Add($noWrapStyleName)
$phoneStyle.IncludeAlignment = $true
$phoneStyle.WrapText = $false
$phoneStyle.Font.Bold = $true
$phoneStyle.Font.Size = 12
$phoneStyle.Font.Name = "Arial"
$phoneStyle.Font.Color = $RED
$ws.Range("H1").Style = $noWrapStyleName
$ws.Range("H1").HorizontalAlignment = $xlCenter
